$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Beta Breeze"
$ws.Range("B3").Value = "Sembahwang"
$ws.Range("C3").Value = "2-Room"
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 300000
$ws.Range("F3").Value = "3-Room"
$ws.Range("G3").Value = 4
$ws.Range("H3").Value = 400000
$ws.Range("I2:J2").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("I3").Value2 = 45703
$ws.Range("J3").Value2 = 45736
$ws.Range("K3").Value = "Michael"
$ws.Range("L3").Value = 4
$ws.Range("M3").Value = "David"

$ws.Range("M4").Select()
